# Fixing the pfc footprints and relating them to their lcsc counterparts.
# Adds a new LCSC part row (row 26) to the BOM sheet and tweaks a handful
# of column widths / view settings to better frame the now-taller sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 26: FlyWin MO3W-82K resistor (C601908) -----------------------
$ws.Range("A26").Value = "C601908"
$ws.Range("B26").Value = "MO3W-82K±5%-QT73"
$ws.Range("C26").Value = "FlyWin"
$ws.Range("D26").Value = "axial"
$ws.Range("F26").Value = "82kΩ ±5% 3W ±350ppm/℃ Axial Metal Oxide Resistors"
$ws.Range("G26").Value = "yes"
$ws.Range("H26").Value = 10
$ws.Range("I26").Value = 10
$ws.Range("J26").Value = 0.0487
$ws.Range("K26").Value = 0.49
$ws.Range("L26").Value = "lcsc.com/product-detail/Metal-Oxide-Resistors_FlyWin-MO3W-82K-5-QT73_C601908.html"

# Row formatting: taller row, LCSC part + unit price wrap, price right aligned
$ws.Rows.Item(26).RowHeight = 15.7
$ws.Range("A26").WrapText = $true
$ws.Range("J26").WrapText = $true
$ws.Range("J26").HorizontalAlignment = -4152

# --- Column width retouch for the now-longer description column ----------
$ws.Columns.Item(3).ColumnWidth = 39.5
$ws.Columns.Item(4).ColumnWidth = 13
$ws.Columns.Item(5).ColumnWidth = 4.166667
$ws.Columns.Item(6).ColumnWidth = 19.5
$ws.Columns.Item(9).ColumnWidth = 9

# --- View: keep only the header frozen, scroll back up, select A12 -------
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$ws.Range("A12").Select()
